# Reran questionnaire analysis with corrected questionnaire data.
$wb = $excel.ActiveWorkbook

$wsNormality = $wb.Worksheets.Item("normality")
$wsPairwise  = $wb.Worksheets.Item("pairwise_tests")

# --- normality sheet: updated W / pval statistics ---
$wsNormality.Range("B3").Value = 0.9817
$wsNormality.Range("C3").Value = 0.9893999999999999
$wsNormality.Range("B4").Value = 0.8204
$wsNormality.Range("C4").Value = 0.012

# --- pairwise_tests sheet: updated T, dof, p-unc, BF10, hedges ---
$wsPairwise.Range("G3").Value = 0.7823
$wsPairwise.Range("H3").Value = 22.1656
$wsPairwise.Range("J3").Value = 0.4423

# K3 (BF10) is stored as text in the workbook ("0.369" -> "0.461").
# Force text entry (Excel would otherwise infer a number) and then
# restore the cell's original (default) style so only the value changes.
$wsPairwise.Range("K3").NumberFormat = "@"
$wsPairwise.Range("K3").Value = "0.461"
$wsPairwise.Range("K3").Style = "Normal"

$wsPairwise.Range("L3").Value = 0.2994
